# Update cryptocurrency price values in column D per the source diff.
# Values are textual (stored as inline strings in the source) so we use an
# apostrophe text-prefix to force text storage, then ClearFormats to drop the
# transient quotePrefix style IronCalc/Excel would otherwise leave behind -
# keeping cell formatting identical to the untouched cells around it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'250.29"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = "'24.12"
$ws.Range("D3").ClearFormats()
$ws.Range("D4").Value = "'6.032"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").Value = "'0.05977"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'3.424"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'6.561"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").Value = "'1.325"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'0.7978"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").Value = "'0.1493"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").Value = "'0.07929"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.03348"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'0.03095"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").Value = "'0.09289"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'3.571"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "'0.001685"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "'0.04766"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "'0.0006104"
$ws.Range("D18").ClearFormats()
$ws.Range("D20").Value = "'0.005696"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").Value = "'0.001072"
$ws.Range("D21").ClearFormats()
$ws.Range("D23").Value = "'3.673"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").Value = "'2.204"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").Value = "'0.3351"
$ws.Range("D25").ClearFormats()
$ws.Range("D27").Value = "'0.0006490"
$ws.Range("D27").ClearFormats()
$ws.Range("D40").Value = "'0.04437"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "'0.007032"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").Value = "'0.003609"
$ws.Range("D42").ClearFormats()
$ws.Range("D44").Value = "'0.009340"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'0.002466"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").Value = "'0.00005899"
$ws.Range("D46").ClearFormats()
$ws.Range("D48").Value = "'0.7017"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Value = "'0.09604"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").Value = "'0.01012"
$ws.Range("D51").ClearFormats()
